$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.004.08'
$ws.Range('E2').Value = '  -0.33%  '

$ws.Range('D3').Value = '2.419.68'
$ws.Range('E3').Value = '  -0.13%  '

$style = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '562.93'
$ws.Range('D5').Style = $style
$ws.Range('E5').Value = '  +0.43%  '

$style = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.32'
$ws.Range('D6').Style = $style
$ws.Range('E6').Value = '  -0.56%  '

$ws.Range('E7').Value = '  -0.04%  '

$style = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.532'
$ws.Range('D8').Style = $style
$ws.Range('E8').Value = '  -0.34%  '

$ws.Range('E9').Value = '  -0.13%  '

$ws.Range('E10').Value = '  -0.76%  '

$ws.Range('E11').Value = '  -4.13%  '

$style = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.350'
$ws.Range('D12').Style = $style
$ws.Range('E12').Value = '  -0.94%  '

$style = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '26.08'
$ws.Range('D13').Style = $style
$ws.Range('E13').Value = '  -0.08%  '

$ws.Range('E14').Value = '  -1.30%  '

$ws.Range('D15').Value = '2.854.40'
$ws.Range('E15').Value = '  +0.02%  '

$ws.Range('D16').Value = '61.878.86'
$ws.Range('E16').Value = '  -0.28%  '

$ws.Range('D17').Value = '2.410.80'
$ws.Range('E17').Value = '  -0.35%  '

$ws.Range('E18').Value = '  +1.22%  '

$style = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '323.16'
$ws.Range('D19').Style = $style
$ws.Range('E19').Value = '  -0.44%  '

$style = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.83'
$ws.Range('D20').Style = $style
$ws.Range('E20').Value = '  +0.74%  '

$style = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.13'
$ws.Range('D21').Style = $style
$ws.Range('E21').Value = '  -1.38%  '

$ws.Range('E22').Value = '  +0.00%  '

$style = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '66.81'
$ws.Range('D23').Style = $style
$ws.Range('E23').Value = '  +1.99%  '

$style = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.73'
$ws.Range('D24').Style = $style
$ws.Range('E24').Value = '  +0.02%  '

$style = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.75'
$ws.Range('D25').Style = $style
$ws.Range('E25').Value = '  -2.52%  '

$style = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '555.45'
$ws.Range('D26').Style = $style
$ws.Range('E26').Value = '  -5.66%  '

$ws.Range('D27').Value = '2.535.83'
$ws.Range('E27').Value = '  -0.28%  '

$ws.Range('E28').Value = '  -0.04%  '

$ws.Range('D29').Value = '0.0₃0936'
$ws.Range('E29').Value = '  -0.57%  '

$style = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.20'
$ws.Range('D30').Style = $style
$ws.Range('E30').Value = '  -0.95%  '

$style = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.39'
$ws.Range('D31').Style = $style
$ws.Range('E31').Value = '  -4.29%  '

$ws.Range('E32').Value = '  -1.60%  '

$style = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.88'
$ws.Range('D33').Style = $style
$ws.Range('E33').Value = '  -0.70%  '

$style = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.51'
$ws.Range('D34').Style = $style
$ws.Range('E34').Value = '  -3.85%  '

$ws.Range('E35').Value = '  -0.05%  '

$style = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.74'
$ws.Range('D36').Style = $style
$ws.Range('E36').Value = '  -0.77%  '

$ws.Range('E37').Value = '  -1.67%  '

$style = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '153.45'
$ws.Range('D38').Style = $style
$ws.Range('E38').Value = '  -0.35%  '

$style = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.42'
$ws.Range('D39').Style = $style
$ws.Range('E39').Value = '  -5.16%  '

$style = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.53'
$ws.Range('D40').Style = $style
$ws.Range('E40').Value = '  -1.07%  '

$style = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.82'
$ws.Range('D41').Style = $style
$ws.Range('E41').Value = '  +0.19%  '

$style = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.999'
$ws.Range('D42').Style = $style
$ws.Range('E42').Value = '  -0.05%  '

$style = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '147.14'
$ws.Range('D43').Style = $style
$ws.Range('E43').Value = '  -2.51%  '

$style = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.24'
$ws.Range('D44').Style = $style
$ws.Range('E44').Value = '  -5.71%  '

$style = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.64'
$ws.Range('D45').Style = $style
$ws.Range('E45').Value = '  -0.38%  '

$style = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0528'
$ws.Range('D46').Style = $style
$ws.Range('E46').Value = '  -2.10%  '

$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$style = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.593'
$ws.Range('D47').Style = $style
$ws.Range('E47').Value = '  +0.02%  '

$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$style = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '19.80'
$ws.Range('D48').Style = $style
$ws.Range('E48').Value = '  -2.67%  '

$style = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0921'
$ws.Range('D49').Style = $style
$ws.Range('E49').Value = '  -0.40%  '

$ws.Range('E50').Value = '  -0.77%  '

$ws.Range('E51').Value = '  +4.59%  '
